# Update dashboards - 2025-10-17
# Pulls a fresh day's reading into the rolling "Present / Lag1-4" window for
# several FRED series on the "Aguilar Prototype" sheet, shifting older
# observations right and dropping the oldest one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-LatestDate {
    param($rowNum, $dateText)
    # Leading apostrophe forces the date-shaped text to stay text (matches
    # the existing cells in column N, which are plain strings, not real
    # Excel dates).
    $ws.Range("N$rowNum").Value = "'" + $dateText
}

# Row 29: 5yr, 5yr Forward (T5YIFR) -- new "Present" reading only
Set-LatestDate 29 "2025-10-16"
$ws.Range("Q29").Value = 2.27

# Row 30: 10yr TIPS (T10YIE) -- new day, values shift right, oldest drops
Set-LatestDate 30 "2025-10-16"
$ws.Range("Q30").Value = 2.28
$ws.Range("R30").Value = 2.29
$ws.Range("S30").Value = 2.3
$ws.Range("T30").Value = 0
$ws.Range("U30").Value = 2.3

# Row 48: 2y UST (DGS2) -- new day, values shift right, oldest drops
Set-LatestDate 48 "2025-10-15"
$ws.Range("Q48").Value = 3.5
$ws.Range("R48").Value = 3.48
$ws.Range("S48").Value = 0
$ws.Range("T48").Value = 3.52
$ws.Range("U48").Value = 3.6

# Row 49: 5y UST (DGS5) -- new day, values shift right, oldest drops
Set-LatestDate 49 "2025-10-15"
$ws.Range("Q49").Value = 3.63
$ws.Range("R49").Value = 3.6
$ws.Range("S49").Value = 0
$ws.Range("T49").Value = 3.65
$ws.Range("U49").Value = 3.74

# Row 50: 10y UST (DGS10) -- new day, values shift right, oldest drops
Set-LatestDate 50 "2025-10-15"
$ws.Range("Q50").Value = 4.05
$ws.Range("R50").Value = 4.03
$ws.Range("S50").Value = 0
$ws.Range("T50").Value = 4.05
$ws.Range("U50").Value = 4.14

# Row 51: 30y Mtg. (MORTGAGE30US) -- new week, values shift right, oldest drops
Set-LatestDate 51 "2025-10-16"
$ws.Range("Q51").Value = 6.27
$ws.Range("R51").Value = 6.3
$ws.Range("S51").Value = 6.34
$ws.Range("T51").Value = 6.3
$ws.Range("U51").Value = 6.26

# Row 52: BAA (DBAA) -- new day, values shift right, oldest drops
Set-LatestDate 52 "2025-10-15"
$ws.Range("Q52").Value = 5.73
$ws.Range("R52").Value = 5.74
$ws.Range("S52").Value = 0
$ws.Range("T52").Value = 5.77
$ws.Range("U52").Value = 5.83
